# Auto-generated Excel COM-interop script
# Updates static market-price / profit columns (H:N) on each job sheet,
# matching a refreshed data pull from the market-board API.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H2").Value = 2302.7334
$ws.Range("I2").Value = 5075
$ws.Range("J2").Value = 1294.6364
$ws.Range("K2").Value = 5075
$ws.Range("L2").Value = 1294.6364
$ws.Range("M2").Value = -4962
$ws.Range("N2").Value = -1520.6364
$ws.Range("H19").Value = 1190
$ws.Range("I19").Value = 987.5
$ws.Range("K19").Value = 987.5
$ws.Range("M19").Value = -812.5
$ws.Range("H32").Value = 6336
$ws.Range("J32").Value = 6336
$ws.Range("L32").Value = 6336
$ws.Range("N32").Value = -6988
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1825
$ws.Range("H43").Value = 5759.4
$ws.Range("I43").Value = 5385
$ws.Range("K43").Value = 5385
$ws.Range("M43").Value = -5316
$ws.Range("H51").Value = 11549.292
$ws.Range("I51").Value = 10300.1
$ws.Range("J51").Value = 12441.571
$ws.Range("K51").Value = 10300.1
$ws.Range("L51").Value = 12441.571
$ws.Range("M51").Value = -9816.1
$ws.Range("N51").Value = -13409.571
$ws.Range("H58").Value = 4900.7617
$ws.Range("I58").Value = 1547.4445
$ws.Range("K58").Value = 4642.333500000001
$ws.Range("M58").Value = -4492.333500000001
$ws.Range("H132").Value = 3316.8333
$ws.Range("I132").Value = 2836.6667
$ws.Range("K132").Value = 8510.000100000001
$ws.Range("M132").Value = -5980.000100000001
$ws.Range("H136").Value = 100000
$ws.Range("J136").Value = 100000
$ws.Range("L136").Value = 100000
$ws.Range("N136").Value = -110200

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H43").Value = 125000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 125000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 125000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -125626
$ws.Range("H45").Value = 2781.1667
$ws.Range("I45").Value = 2437.4
$ws.Range("K45").Value = 2437.4
$ws.Range("M45").Value = -2060.4
$ws.Range("H63").Value = 7263.533
$ws.Range("J63").Value = 9149.9
$ws.Range("L63").Value = 9149.9
$ws.Range("N63").Value = -10521.9
$ws.Range("H66").Value = 7263.533
$ws.Range("J66").Value = 9149.9
$ws.Range("L66").Value = 45749.5
$ws.Range("N66").Value = -52613.5
$ws.Range("H110").Value = 5129794.5
$ws.Range("I110").Value = 10990167
$ws.Range("K110").Value = 10990167
$ws.Range("M110").Value = -10988122
$ws.Range("H122").Value = 11056438
$ws.Range("I122").Value = 15003381
$ws.Range("K122").Value = 45010143
$ws.Range("M122").Value = -45007693
$ws.Range("H132").Value = 1621.9231
$ws.Range("J132").Value = 1770.125
$ws.Range("L132").Value = 5310.375
$ws.Range("N132").Value = -10370.375

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H23").Value = 4749.5
$ws.Range("I23").Value = 4500
$ws.Range("J23").Value = 4999
$ws.Range("K23").Value = 4500
$ws.Range("L23").Value = 4999
$ws.Range("M23").Value = -4217
$ws.Range("N23").Value = -5565
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H25").Value = 1970.5714
$ws.Range("J25").Value = 2114.4614
$ws.Range("L25").Value = 2114.4614
$ws.Range("N25").Value = -2462.4614
$ws.Range("H58").Value = 3748.037
$ws.Range("I58").Value = 3386.4707
$ws.Range("K58").Value = 3386.4707
$ws.Range("M58").Value = -3183.4707
$ws.Range("H69").Value = 19999.5
$ws.Range("I69").Value = 19999.5
$ws.Range("K69").Value = 19999.5
$ws.Range("M69").Value = -19250.5
$ws.Range("H72").Value = 19999.5
$ws.Range("I72").Value = 19999.5
$ws.Range("K72").Value = 59998.5
$ws.Range("M72").Value = -56254.5
$ws.Range("H86").Value = 4893.5
$ws.Range("I86").Value = 4893.5
$ws.Range("K86").Value = 4893.5
$ws.Range("M86").Value = -3770.5
$ws.Range("H89").Value = 4893.5
$ws.Range("I89").Value = 4893.5
$ws.Range("K89").Value = 24467.5
$ws.Range("M89").Value = -18851.5
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 5000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -7996
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 5102.3076
$ws.Range("I132").Value = 4853
$ws.Range("K132").Value = 14559
$ws.Range("M132").Value = -12029
$ws.Range("H136").Value = 3748.037
$ws.Range("I136").Value = 3386.4707
$ws.Range("K136").Value = 10159.4121
$ws.Range("M136").Value = -7609.4121

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H12").Value = 190.125
$ws.Range("I12").Value = 375.33334
$ws.Range("J12").Value = 79
$ws.Range("K12").Value = 1126.00002
$ws.Range("L12").Value = 237
$ws.Range("M12").Value = -953.0000199999999
$ws.Range("N12").Value = -583
$ws.Range("H33").Value = 563.5263
$ws.Range("I33").Value = 594.2778
$ws.Range("K33").Value = 3565.6668
$ws.Range("M33").Value = -3282.6668
$ws.Range("H40").Value = 51.090908
$ws.Range("I40").Value = 21.5
$ws.Range("K40").Value = 86
$ws.Range("M40").Value = -17
$ws.Range("H98").Value = 1503.8572
$ws.Range("I98").Value = 2005.6
$ws.Range("J98").Value = 249.5
$ws.Range("K98").Value = 6016.799999999999
$ws.Range("L98").Value = 748.5
$ws.Range("M98").Value = -4518.799999999999
$ws.Range("N98").Value = -3744.5
$ws.Range("H109").Value = 1623.75
$ws.Range("I109").Value = 1665
$ws.Range("J109").Value = 1500
$ws.Range("K109").Value = 4995
$ws.Range("L109").Value = 4500
$ws.Range("M109").Value = -3955
$ws.Range("N109").Value = -6580

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H80").Value = 14999
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 14999
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H113").Value = 4833.3335
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4833.3335
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4833.3335
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9173.333500000001
$ws.Range("H114").Value = 35500
$ws.Range("J114").Value = 35500
$ws.Range("L114").Value = 35500
$ws.Range("N114").Value = -44178
$ws.Range("H122").Value = 93717.55
$ws.Range("I122").Value = 2799.25
$ws.Range("J122").Value = 336166.34
$ws.Range("K122").Value = 8397.75
$ws.Range("L122").Value = 1008499.02
$ws.Range("M122").Value = -5947.75
$ws.Range("N122").Value = -1013399.02
$ws.Range("H132").Value = 2508.742
$ws.Range("I132").Value = 2145.2693
$ws.Range("K132").Value = 6435.8079
$ws.Range("M132").Value = -3905.8079

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 1115.091
$ws.Range("I22").Value = 754
$ws.Range("K22").Value = 754
$ws.Range("M22").Value = -459
$ws.Range("H27").Value = 1115.091
$ws.Range("I27").Value = 754
$ws.Range("K27").Value = 754
$ws.Range("M27").Value = -647
$ws.Range("H30").Value = 2327.875
$ws.Range("J30").Value = 5353.3335
$ws.Range("L30").Value = 5353.3335
$ws.Range("N30").Value = -5569.3335
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4828
$ws.Range("H35").Value = 2444.0908
$ws.Range("I35").Value = 707.1429000000001
$ws.Range("J35").Value = 5483.75
$ws.Range("K35").Value = 707.1429000000001
$ws.Range("L35").Value = 5483.75
$ws.Range("M35").Value = -371.1429000000001
$ws.Range("N35").Value = -6155.75
$ws.Range("H122").Value = 3836
$ws.Range("I122").Value = 2920
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 8760
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -6310
$ws.Range("N122").Value = -27400
$ws.Range("H136").Value = 4005
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4005
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 12015
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -17115

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H8").Value = 13533
$ws.Range("J8").Value = 13533
$ws.Range("L8").Value = 13533
$ws.Range("N8").Value = -13813
$ws.Range("H9").Value = 20000
$ws.Range("I9").Value = 20000
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -19860
$ws.Range("N9").ClearContents()
$ws.Range("H68").Value = 50271
$ws.Range("J68").Value = 50271
$ws.Range("L68").Value = 50271
$ws.Range("N68").Value = -51893
$ws.Range("H71").Value = 50271
$ws.Range("J71").Value = 50271
$ws.Range("L71").Value = 150813
$ws.Range("N71").Value = -158925
$ws.Range("H132").Value = 2002.88
$ws.Range("I132").Value = 1866.9474
$ws.Range("K132").Value = 5600.8422
$ws.Range("M132").Value = -3070.8422
